$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author reworked row 6 (previously a plain numeric AHV/NNSS number with
# no name attached) into a full record, and appended two more rows covering
# additional "ahv_format" edge cases (letters embedded in the AHV number)
# that the pseudonymizer needs to detect.
#
# Cells are filled in the same order the original author typed them, so the
# shared-strings table comes out in the same sequence as the authored file.

$ws.Range("A7").Value = "756AB00000004"
$ws.Range("A8").Value = "7230000XYZ"

$ws.Range("B6").Value = "Dachs"
$ws.Range("C6").Value = "Dario"

$ws.Range("B7").Value = "Eris"
$ws.Range("C7").Value = "Elsa"

$ws.Range("C8").Value = "Fiona"
$ws.Range("B8").Value = "Fichter"

# A6 becomes a text value (it holds a malformed/otherwise-marked AHV number)
$ws.Range("A6").Value = "756.0000.000.004"

# Row 7's AHV cell (A7) should carry the same "NNSS" number format/style as
# the existing A2:A6 entries above it, so copy that formatting across.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Mirror the author's final cursor position in the saved file.
$ws.Range("A6").Select() | Out-Null
